$d = $word.ActiveDocument

$pairs = @(
    @("35×97=3395", "84×37=3108"),
    @("51×95=4845", "61×74=4514"),
    @("29×97=2813", "49×26=1274"),
    @("45×50=2250", "71×88=6248"),
    @("24×60=1440", "54×26=1404"),
    @("47×29=1363", "97×23=2231"),
    @("25×95=2375", "18×77=1386"),
    @("45×66=2970", "82×99=8118"),
    @("78×41=3198", "22×68=1496"),
    @("21×11=231",  "93×30=2790"),
    @("98×99=9702", "99×66=6534"),
    @("16×75=1200", "12×85=1020"),
    @("94×67=6298", "77×78=6006"),
    @("33×26=858",  "27×59=1593"),
    @("67×76=5092", "64×16=1024"),
    @("29×16=464",  "73×75=5475"),
    @("55×82=4510", "22×45=990"),
    @("81×22=1782", "38×77=2926"),
    @("55×34=1870", "76×35=2660"),
    @("11×44=484",  "74×47=3478"),
    @("99×73=7227", "21×34=714"),
    @("18×27=486",  "37×55=2035"),
    @("14×12=168",  "31×69=2139"),
    @("51×20=1020", "79×93=7347"),
    @("79×72=5688", "26×19=494")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
